$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 4-6 with new data
$ws.Range("A4").Value = "March 14th, 2025"
$ws.Range("B4").Value = "Cooking food"

$ws.Range("A5").Value = "March 30th, 2025"
$ws.Range("B5").Value = "Get a new phone"

$ws.Range("A6").Value = "Tomorrow"
$ws.Range("B6").Value = "CSC333"

# Remove rows 7-9 entirely (delete the rows, shifting cells up)
$ws.Range("A7:D9").Delete()
